# Weekly update: insert a new week of "Coliflor" (Mercado Mayorista Lo
# Valledor de Santiago) records at the top of the data block and push the
# rest of the rows down by two (R592 -> R594).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 561; everything that used to live in
# rows 561:592 shifts down to 563:594 (dimension becomes A1:R594).
$ws.Rows("561:562").Insert()

$newDate = Get-Date -Year 2022 -Month 1 -Day 24 -Hour 0 -Minute 0 -Second 0

# New row 561: Coliflor, Primera
$ws.Range("A561").Value = 6
$ws.Range("B561").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C561").Value = 'Metropolitana'
$ws.Range("D561").Value = $newDate
$ws.Range("E561").Value = 13
$ws.Range("F561").Value = 100112008
$ws.Range("G561").Value = 'Coliflor'
$ws.Range("H561").Value = 'Sin especificar'
$ws.Range("I561").Value = 'Primera'
$ws.Range("J561").Value = 5800
$ws.Range("K561").Value = 800
$ws.Range("L561").Value = 900
$ws.Range("M561").Value = 866
$ws.Range("N561").Value = '$/unidad'
$ws.Range("O561").Value = 'Región Metropolitana'
$ws.Range("P561").Value = 866
$ws.Range("Q561").Value = 1
$ws.Range("R561").Value = 'Hortaliza'

# New row 562: Coliflor, Segunda
$ws.Range("A562").Value = 6
$ws.Range("B562").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C562").Value = 'Metropolitana'
$ws.Range("D562").Value = $newDate
$ws.Range("E562").Value = 13
$ws.Range("F562").Value = 100112008
$ws.Range("G562").Value = 'Coliflor'
$ws.Range("H562").Value = 'Sin especificar'
$ws.Range("I562").Value = 'Segunda'
$ws.Range("J562").Value = 1800
$ws.Range("K562").Value = 700
$ws.Range("L562").Value = 700
$ws.Range("M562").Value = 700
$ws.Range("N562").Value = '$/unidad'
$ws.Range("O562").Value = 'Región Metropolitana'
$ws.Range("P562").Value = 700
$ws.Range("Q562").Value = 1
$ws.Range("R562").Value = 'Hortaliza'
